$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 9).Value = 0.9179281773574478
$ws.Cells.Item(2, 10).Value = 0.9179281773574478
$ws.Cells.Item(2, 13).Value = 1.994783
$ws.Cells.Item(2, 14).Value = 5.984349
$ws.Cells.Item(2, 15).Value = 0.2551127970612924
$ws.Cells.Item(2, 16).Value = 0.2551127970612924
$ws.Cells.Item(2, 17).Value = 1.230337604246333
$ws.Cells.Item(2, 18).Value = 11.073038438217
$ws.Cells.Item(2, 19).Value = 0.2341752248270326
$ws.Cells.Item(2, 20).Value = 0.2341752248270326

# Row 3
$ws.Cells.Item(3, 9).Value = 0.9179281773574478
$ws.Cells.Item(3, 10).Value = 0.9179281773574478
$ws.Cells.Item(3, 13).Value = 4.293075666666667
$ws.Cells.Item(3, 15).Value = 0.5490414452695386
$ws.Cells.Item(3, 16).Value = 0.5490414452695387
$ws.Cells.Item(3, 19).Value = 0.5039806131499664
$ws.Cells.Item(3, 20).Value = 0.5039806131499666

# Row 4
$ws.Cells.Item(4, 9).Value = 0.9179281773574478
$ws.Cells.Item(4, 10).Value = 0.9179281773574478
$ws.Cells.Item(4, 13).Value = 1.506399333333333
$ws.Cells.Item(4, 14).Value = 4.519197999999999
$ws.Cells.Item(4, 15).Value = 0.1926534101292887
$ws.Cells.Item(4, 16).Value = 0.1926534101292887
$ws.Cells.Item(4, 17).Value = 0.9291134658815555
$ws.Cells.Item(4, 18).Value = 8.362021192933998
$ws.Cells.Item(4, 19).Value = 0.1768419936216748
$ws.Cells.Item(4, 20).Value = 0.1768419936216748

# Row 5
$ws.Cells.Item(5, 9).Value = 0.9179281773574478
$ws.Cells.Item(5, 10).Value = 0.9179281773574478
$ws.Cells.Item(5, 13).Value = 0.02496166666666666
$ws.Cells.Item(5, 14).Value = 0.07488499999999999
$ws.Cells.Item(5, 15).Value = 0.003192347539880258
$ws.Cells.Item(5, 16).Value = 0.003192347539880258
$ws.Cells.Item(5, 17).Value = 0.01539579852277778
$ws.Cells.Item(5, 18).Value = 0.138562186705
$ws.Cells.Item(5, 19).Value = 0.002930345758773817
$ws.Cells.Item(5, 20).Value = 0.002930345758773818

# Row 6
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.055146
$ws.Cells.Item(6, 8).Value = 0.165438
$ws.Cells.Item(6, 9).Value = 0.08207182264255215
$ws.Cells.Item(6, 10).Value = 0.08207182264255215
$ws.Cells.Item(6, 13).Value = 1.994783
$ws.Cells.Item(6, 14).Value = 5.984349
$ws.Cells.Item(6, 15).Value = 0.2551127970612924
$ws.Cells.Item(6, 16).Value = 0.2551127970612924
$ws.Cells.Item(6, 17).Value = 0.110004303318
$ws.Cells.Item(6, 18).Value = 0.990038729862
$ws.Cells.Item(6, 19).Value = 0.02093757223425979
$ws.Cells.Item(6, 20).Value = 0.02093757223425979

# Row 7
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.055146
$ws.Cells.Item(7, 8).Value = 0.165438
$ws.Cells.Item(7, 9).Value = 0.08207182264255215
$ws.Cells.Item(7, 10).Value = 0.08207182264255215
$ws.Cells.Item(7, 13).Value = 4.293075666666667
$ws.Cells.Item(7, 15).Value = 0.5490414452695386
$ws.Cells.Item(7, 16).Value = 0.5490414452695387
$ws.Cells.Item(7, 17).Value = 0.236745950714
$ws.Cells.Item(7, 18).Value = 2.130713556426
$ws.Cells.Item(7, 19).Value = 0.04506083211957208
$ws.Cells.Item(7, 20).Value = 0.04506083211957208

# Row 8
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.055146
$ws.Cells.Item(8, 8).Value = 0.165438
$ws.Cells.Item(8, 9).Value = 0.08207182264255215
$ws.Cells.Item(8, 10).Value = 0.08207182264255215
$ws.Cells.Item(8, 13).Value = 1.506399333333333
$ws.Cells.Item(8, 14).Value = 4.519197999999999
$ws.Cells.Item(8, 15).Value = 0.1926534101292887
$ws.Cells.Item(8, 16).Value = 0.1926534101292887
$ws.Cells.Item(8, 17).Value = 0.08307189763599999
$ws.Cells.Item(8, 18).Value = 0.7476470787239999
$ws.Cells.Item(8, 19).Value = 0.01581141650761384
$ws.Cells.Item(8, 20).Value = 0.01581141650761384

# Row 9
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.055146
$ws.Cells.Item(9, 8).Value = 0.165438
$ws.Cells.Item(9, 9).Value = 0.08207182264255215
$ws.Cells.Item(9, 10).Value = 0.08207182264255215
$ws.Cells.Item(9, 13).Value = 0.02496166666666666
$ws.Cells.Item(9, 14).Value = 0.07488499999999999
$ws.Cells.Item(9, 15).Value = 0.003192347539880258
$ws.Cells.Item(9, 16).Value = 0.003192347539880258
$ws.Cells.Item(9, 17).Value = 0.00137653607
$ws.Cells.Item(9, 18).Value = 0.01238882463
$ws.Cells.Item(9, 19).Value = 0.0002620017811064402
$ws.Cells.Item(9, 20).Value = 0.0002620017811064402
